$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 4.33
$ws.Range("K2").Value = 2.38
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = 1.62
$ws.Range("R2").Value = 2.3
$ws.Range("S2").Value = 2.04
$ws.Range("T2").Value = 1.86
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 1.3
$ws.Range("X2").Value = 3.4
$ws.Range("Y2").Value = 1.57
$ws.Range("Z2").Value = 2.25
$ws.Range("AI2").Value = 15
$ws.Range("AM2").Value = 11
$ws.Range("AO2").Value = 17
$ws.Range("AQ2").Value = 23

# Row 3
$ws.Range("M3").Value = 1.05
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 2.62
$ws.Range("V3").Value = 1.15

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.5
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.38
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.37
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.67
$ws.Range("U4").Value = 4
$ws.Range("V4").Value = 1.22
$ws.Range("W4").Value = 1.5
$ws.Range("X4").Value = 2.5
$ws.Range("Y4").Value = 2.1
$ws.Range("Z4").Value = 1.67
$ws.Range("AA4").Value = 6
$ws.Range("AB4").Value = 7
$ws.Range("AD4").Value = 13
$ws.Range("AF4").Value = 34
$ws.Range("AG4").Value = 8
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 67
$ws.Range("AK4").Value = 501
$ws.Range("AL4").Value = 11
$ws.Range("AN4").Value = 17
$ws.Range("AQ4").Value = 51

# Row 5
$ws.Range("G5").Value = 4.75
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 1.65
$ws.Range("J5").Value = 5
$ws.Range("L5").Value = 2.25
$ws.Range("M5").Value = 1.03
$ws.Range("O5").Value = 1.19
$ws.Range("P5").Value = 4
$ws.Range("V5").Value = 1.33
$ws.Range("AC5").Value = 15
$ws.Range("AI5").Value = 15
$ws.Range("AK5").Value = 201
$ws.Range("AO5").Value = 13

# Row 7
$ws.Range("G7").Value = 1.62
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 5.5
$ws.Range("J7").Value = 2.2
$ws.Range("K7").Value = 2.3
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 12
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4
$ws.Range("Q7").Value = 1.75
$ws.Range("R7").Value = 2.05
$ws.Range("U7").Value = 2.75
$ws.Range("V7").Value = 1.4
$ws.Range("Y7").Value = 1.8
$ws.Range("Z7").Value = 1.95
$ws.Range("AA7").Value = 7.5
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 12
$ws.Range("AE7").Value = 13
$ws.Range("AG7").Value = 12
$ws.Range("AH7").Value = 7.5
$ws.Range("AL7").Value = 15
